$d = $word.ActiveDocument

# Remove the now-unused "Abstract Title" custom paragraph style entirely.
$abstractTitleStyle = $d.Styles.Item("Abstract Title")
$abstractTitleStyle.Delete()

# Update the "Abstract" style's paragraph spacing: before 5pt -> 15pt
# (w:before 100 -> 300, i.e. twentieths-of-a-point / 20 = points).
$abstractStyle = $d.Styles.Item("Abstract")
$abstractStyle.ParagraphFormat.SpaceBefore = 15
